# Applies Rachel's edits to the logistic-regression coefficient / expression-group table:
#  - renames the "val_logitlasso" column header to "val_lasso"
#  - refreshes the var/coef/group/wordCount data for every row (rows 2-106)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'var'
$ws.Range("B1").Value = 'val_lasso'
$ws.Range("C1").Value = 'DoD_Dicts'
$ws.Range("D1").Value = 'wordCount'
$ws.Range("A2").Value = 'personal'
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = '{''Individualism''}'
$ws.Range("D2").Value = 96
$ws.Range("A3").Value = 'unavoidable'
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = '{''Doomism''}'
$ws.Range("D3").Value = 6
$ws.Range("A4").Value = 'breakthrough'
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = '{''Tech_optimism''}'
$ws.Range("D4").Value = 53
$ws.Range("A5").Value = 'overuse'
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = '{''Free_rider''}'
$ws.Range("D5").Value = 18
$ws.Range("A6").Value = 'invent'
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = '{''Tech_optimism''}'
$ws.Range("D6").Value = 5
$ws.Range("A7").Value = 'disproportionate'
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = '{''Social_justice''}'
$ws.Range("D7").Value = 21
$ws.Range("A8").Value = 'volunteer'
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = '{''Carrots''}'
$ws.Range("D8").Value = 304
$ws.Range("A9").Value = 'uncertainty'
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = '{''Doomism''}'
$ws.Range("D9").Value = 17
$ws.Range("A10").Value = 'living standards'
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = '{''Well_being''}'
$ws.Range("D10").Value = 3
$ws.Range("A11").Value = 'marketplace'
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = '{''Carrots''}'
$ws.Range("D11").Value = 40
$ws.Range("A12").Value = 'catastrophe'
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = '{''Doomism''}'
$ws.Range("D12").Value = 31
$ws.Range("A13").Value = 'mutually beneficial'
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = '{''Carrots''}'
$ws.Range("D13").Value = 6
$ws.Range("A14").Value = 'infeasible'
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = '{''Change_impossible''}'
$ws.Range("D14").Value = 3
$ws.Range("A15").Value = 'low income'
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = '{''Social_justice''}'
$ws.Range("D15").Value = 11
$ws.Range("A16").Value = 'compete'
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = '{''Whataboutism''}'
$ws.Range("D16").Value = 32
$ws.Range("A17").Value = 'socioeconomic'
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = '{''Social_justice''}'
$ws.Range("D17").Value = 18
$ws.Range("A18").Value = 'compromise'
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = '{''Perfect_policy''}'
$ws.Range("D18").Value = 31
$ws.Range("A19").Value = 'long term'
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = '{''Talk_no_action''}'
$ws.Range("D19").Value = 67
$ws.Range("A20").Value = 'overused'
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = '{''Free_rider''}'
$ws.Range("D20").Value = 2
$ws.Range("A21").Value = 'lost'
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = '{''Well_being''}'
$ws.Range("D21").Value = 194
$ws.Range("A22").Value = 'carbon footprint'
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = '{''Whataboutism''}'
$ws.Range("D22").Value = 21
$ws.Range("A23").Value = 'energy needs'
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = '{''Well_being''}'
$ws.Range("D23").Value = 36
$ws.Range("A24").Value = 'sequestration'
$ws.Range("B24").Value = 0
$ws.Range("C24").Value = '{''Tech_optimism''}'
$ws.Range("D24").Value = 27
$ws.Range("A25").Value = 'fear'
$ws.Range("B25").Value = 0
$ws.Range("C25").Value = '{''Doomism''}'
$ws.Range("D25").Value = 32
$ws.Range("A26").Value = 'commitment'
$ws.Range("B26").Value = 0
$ws.Range("C26").Value = '{''Talk_no_action''}'
$ws.Range("D26").Value = 629
$ws.Range("A27").Value = 'reliable energy'
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = '{''FF_solutionism''}'
$ws.Range("D27").Value = 1
$ws.Range("A28").Value = 'disruption'
$ws.Range("B28").Value = 0
$ws.Range("C28").Value = '{''Social_justice''}'
$ws.Range("D28").Value = 51
$ws.Range("A29").Value = 'overburdened'
$ws.Range("B29").Value = 0
$ws.Range("C29").Value = '{''Carrots'', ''Free_rider''}'
$ws.Range("D29").Value = 6
$ws.Range("A30").Value = 'renewable natural'
$ws.Range("B30").Value = 0
$ws.Range("C30").Value = '{''FF_solutionism''}'
$ws.Range("D30").Value = 2
$ws.Range("A31").Value = 'footprint'
$ws.Range("B31").Value = 0
$ws.Range("C31").Value = '{''Individualism''}'
$ws.Range("D31").Value = 75
$ws.Range("A32").Value = 'promised'
$ws.Range("B32").Value = 0
$ws.Range("C32").Value = '{''Talk_no_action''}'
$ws.Range("D32").Value = 24
$ws.Range("A33").Value = 'private sector'
$ws.Range("B33").Value = 0
$ws.Range("C33").Value = '{''Carrots''}'
$ws.Range("D33").Value = 75
$ws.Range("A34").Value = 'take advantage'
$ws.Range("B34").Value = 0
$ws.Range("C34").Value = '{''Free_rider''}'
$ws.Range("D34").Value = 25
$ws.Range("A35").Value = 'cautious approach'
$ws.Range("B35").Value = 0
$ws.Range("C35").Value = '{''Perfect_policy''}'
$ws.Range("D35").Value = 2
$ws.Range("A36").Value = 'incentive'
$ws.Range("B36").Value = 0
$ws.Range("C36").Value = '{''Carrots''}'
$ws.Range("D36").Value = 57
$ws.Range("A37").Value = 'creator'
$ws.Range("B37").Value = 0
$ws.Range("C37").Value = '{''Doomism''}'
$ws.Range("D37").Value = 7
$ws.Range("A38").Value = 'carbon intensity'
$ws.Range("B38").Value = 0
$ws.Range("C38").Value = '{''FF_solutionism''}'
$ws.Range("D38").Value = 3
$ws.Range("A39").Value = 'energy management'
$ws.Range("B39").Value = 0
$ws.Range("C39").Value = '{''Tech_optimism''}'
$ws.Range("D39").Value = 8
$ws.Range("A40").Value = 'innovation'
$ws.Range("B40").Value = 0
$ws.Range("C40").Value = '{''Tech_optimism''}'
$ws.Range("D40").Value = 188
$ws.Range("A41").Value = 'prescribe'
$ws.Range("B41").Value = 0
$ws.Range("C41").Value = '{''Carrots''}'
$ws.Range("D41").Value = 2
$ws.Range("A42").Value = 'inflation'
$ws.Range("B42").Value = 0
$ws.Range("C42").Value = '{''Social_justice''}'
$ws.Range("D42").Value = 2
$ws.Range("A43").Value = 'poor'
$ws.Range("B43").Value = 0
$ws.Range("C43").Value = '{''Social_justice''}'
$ws.Range("D43").Value = 98
$ws.Range("A44").Value = 'affordable'
$ws.Range("B44").Value = 0
$ws.Range("C44").Value = '{''Social_justice''}'
$ws.Range("D44").Value = 70
$ws.Range("A45").Value = 'incentives'
$ws.Range("B45").Value = 0
$ws.Range("C45").Value = '{''Carrots''}'
$ws.Range("D45").Value = 227
$ws.Range("A46").Value = 'ambition'
$ws.Range("B46").Value = 0
$ws.Range("C46").Value = '{''Talk_no_action''}'
$ws.Range("D46").Value = 9
$ws.Range("A47").Value = 'lowincome'
$ws.Range("B47").Value = 0
$ws.Range("C47").Value = '{''Social_justice''}'
$ws.Range("D47").Value = 101
$ws.Range("A48").Value = 'exploit'
$ws.Range("B48").Value = 0
$ws.Range("C48").Value = '{''Free_rider''}'
$ws.Range("D48").Value = 22
$ws.Range("A49").Value = 'grim'
$ws.Range("B49").Value = 0
$ws.Range("C49").Value = '{''Doomism''}'
$ws.Range("D49").Value = 17
$ws.Range("A50").Value = 'bipartisan'
$ws.Range("B50").Value = 0
$ws.Range("C50").Value = '{''Perfect_policy''}'
$ws.Range("D50").Value = 107
$ws.Range("A51").Value = 'total emissions'
$ws.Range("B51").Value = 0
$ws.Range("C51").Value = '{''Whataboutism''}'
$ws.Range("D51").Value = 7
$ws.Range("A52").Value = 'regressive'
$ws.Range("B52").Value = 0
$ws.Range("C52").Value = '{''Social_justice''}'
$ws.Range("D52").Value = 4
$ws.Range("A53").Value = 'negligible'
$ws.Range("B53").Value = 0
$ws.Range("C53").Value = '{''Whataboutism''}'
$ws.Range("D53").Value = 4
$ws.Range("A54").Value = 'natural gas'
$ws.Range("B54").Value = 0
$ws.Range("C54").Value = '{''FF_solutionism''}'
$ws.Range("D54").Value = 105
$ws.Range("A55").Value = 'consequence'
$ws.Range("B55").Value = 0
$ws.Range("C55").Value = '{''Well_being''}'
$ws.Range("D55").Value = 14
$ws.Range("A56").Value = 'costly'
$ws.Range("B56").Value = 0
$ws.Range("C56").Value = '{''Social_justice''}'
$ws.Range("D56").Value = 49
$ws.Range("A57").Value = 'failure'
$ws.Range("B57").Value = 0
$ws.Range("C57").Value = '{''Change_impossible''}'
$ws.Range("D57").Value = 82
$ws.Range("A58").Value = 'extreme'
$ws.Range("B58").Value = 0
$ws.Range("C58").Value = '{''Doomism''}'
$ws.Range("D58").Value = 94
$ws.Range("A59").Value = 'impossible'
$ws.Range("B59").Value = 0
$ws.Range("C59").Value = '{''Doomism''}'
$ws.Range("D59").Value = 41
$ws.Range("A60").Value = 'tragedy commons'
$ws.Range("B60").Value = 0
$ws.Range("C60").Value = '{''Free_rider''}'
$ws.Range("D60").Value = 2
$ws.Range("A61").Value = 'unfair'
$ws.Range("B61").Value = 0
$ws.Range("C61").Value = '{''Social_justice''}'
$ws.Range("D61").Value = 5
$ws.Range("A62").Value = 'nonpartisan'
$ws.Range("B62").Value = 0
$ws.Range("C62").Value = '{''Perfect_policy''}'
$ws.Range("D62").Value = 23
$ws.Range("A63").Value = 'target'
$ws.Range("B63").Value = 0
$ws.Range("C63").Value = '{''Whataboutism''}'
$ws.Range("D63").Value = 136
$ws.Range("A64").Value = 'strife'
$ws.Range("B64").Value = 0
$ws.Range("C64").Value = '{''Social_justice''}'
$ws.Range("D64").Value = 3
$ws.Range("A65").Value = 'fusion'
$ws.Range("B65").Value = 0
$ws.Range("C65").Value = '{''Tech_optimism''}'
$ws.Range("D65").Value = 3
$ws.Range("A66").Value = 'best interest'
$ws.Range("B66").Value = 0
$ws.Range("C66").Value = '{''Social_justice'', ''Well_being''}'
$ws.Range("D66").Value = 3
$ws.Range("A67").Value = 'sacrifice'
$ws.Range("B67").Value = 0
$ws.Range("C67").Value = '{''Individualism''}'
$ws.Range("D67").Value = 13
$ws.Range("A68").Value = 'committed'
$ws.Range("B68").Value = 0
$ws.Range("C68").Value = '{''Talk_no_action''}'
$ws.Range("D68").Value = 393
$ws.Range("A69").Value = 'promises'
$ws.Range("B69").Value = 0
$ws.Range("C69").Value = '{''Talk_no_action''}'
$ws.Range("D69").Value = 39
$ws.Range("A70").Value = 'burden'
$ws.Range("B70").Value = 0
$ws.Range("C70").Value = '{''Social_justice'', ''Carrots'', ''Free_rider''}'
$ws.Range("D70").Value = 35
$ws.Range("A71").Value = 'hinder'
$ws.Range("B71").Value = 0
$ws.Range("C71").Value = '{''Social_justice''}'
$ws.Range("D71").Value = 3
$ws.Range("A72").Value = 'vulnerable'
$ws.Range("B72").Value = 0
$ws.Range("C72").Value = '{''Well_being''}'
$ws.Range("D72").Value = 193
$ws.Range("A73").Value = 'fate'
$ws.Range("B73").Value = 0
$ws.Range("C73").Value = '{''Doomism''}'
$ws.Range("D73").Value = 44
$ws.Range("A74").Value = 'exploiting'
$ws.Range("B74").Value = 0
$ws.Range("C74").Value = '{''Free_rider''}'
$ws.Range("D74").Value = 11
$ws.Range("A75").Value = 'commitments'
$ws.Range("B75").Value = 0
$ws.Range("C75").Value = '{''Talk_no_action''}'
$ws.Range("D75").Value = 190
$ws.Range("A76").Value = 'irreversible'
$ws.Range("B76").Value = 0
$ws.Range("C76").Value = '{''Doomism''}'
$ws.Range("D76").Value = 17
$ws.Range("A77").Value = 'appetite'
$ws.Range("B77").Value = 0
$ws.Range("C77").Value = '{''Whataboutism''}'
$ws.Range("D77").Value = 8
$ws.Range("A78").Value = 'unimaginable'
$ws.Range("B78").Value = 0
$ws.Range("C78").Value = '{''Change_impossible''}'
$ws.Range("D78").Value = 4
$ws.Range("A79").Value = 'threat'
$ws.Range("B79").Value = 0
$ws.Range("C79").Value = '{''Well_being''}'
$ws.Range("D79").Value = 367
$ws.Range("A80").Value = 'consumer choice'
$ws.Range("B80").Value = 0
$ws.Range("C80").Value = '{''Individualism''}'
$ws.Range("D80").Value = 2
$ws.Range("A81").Value = 'promise'
$ws.Range("B81").Value = 0
$ws.Range("C81").Value = '{''Talk_no_action''}'
$ws.Range("D81").Value = 93
$ws.Range("A82").Value = 'voluntary'
$ws.Range("B82").Value = 0
$ws.Range("C82").Value = '{''Carrots''}'
$ws.Range("D82").Value = 79
$ws.Range("A83").Value = 'adapt'
$ws.Range("B83").Value = 0
$ws.Range("C83").Value = '{''Doomism''}'
$ws.Range("D83").Value = 80
$ws.Range("A84").Value = 'shared'
$ws.Range("B84").Value = 0
$ws.Range("C84").Value = '{''Free_rider''}'
$ws.Range("D84").Value = 229
$ws.Range("A85").Value = 'insurmountable'
$ws.Range("B85").Value = 0
$ws.Range("C85").Value = '{''Change_impossible''}'
$ws.Range("D85").Value = 4
$ws.Range("A86").Value = 'exploited'
$ws.Range("B86").Value = 0
$ws.Range("C86").Value = '{''Free_rider''}'
$ws.Range("D86").Value = 16
$ws.Range("A87").Value = 'research development'
$ws.Range("B87").Value = 0
$ws.Range("C87").Value = '{''Tech_optimism''}'
$ws.Range("D87").Value = 9
$ws.Range("A88").Value = 'economic prosperity'
$ws.Range("B88").Value = 0
$ws.Range("C88").Value = '{''Well_being''}'
$ws.Range("D88").Value = 6
$ws.Range("A89").Value = 'cleaner fuels'
$ws.Range("B89").Value = 0
$ws.Range("C89").Value = '{''FF_solutionism''}'
$ws.Range("D89").Value = 8
$ws.Range("A90").Value = 'low carbon'
$ws.Range("B90").Value = 0
$ws.Range("C90").Value = '{''FF_solutionism''}'
$ws.Range("D90").Value = 6
$ws.Range("A91").Value = 'share'
$ws.Range("B91").Value = 0
$ws.Range("C91").Value = '{''Free_rider''}'
$ws.Range("D91").Value = 396
$ws.Range("A92").Value = 'unequivocal'
$ws.Range("B92").Value = 0
$ws.Range("C92").Value = '{''Doomism''}'
$ws.Range("D92").Value = 6
$ws.Range("A93").Value = 'investment'
$ws.Range("B93").Value = 0
$ws.Range("C93").Value = '{''Tech_optimism''}'
$ws.Range("D93").Value = 440
$ws.Range("A94").Value = 'disruptive'
$ws.Range("B94").Value = 0
$ws.Range("C94").Value = '{''Well_being''}'
$ws.Range("D94").Value = 11
$ws.Range("A95").Value = 'imminent'
$ws.Range("B95").Value = 0
$ws.Range("C95").Value = '{''Tech_optimism''}'
$ws.Range("D95").Value = 22
$ws.Range("A96").Value = 'resignation'
$ws.Range("B96").Value = 0
$ws.Range("C96").Value = '{''Doomism''}'
$ws.Range("D96").Value = 4
$ws.Range("A97").Value = 'individual'
$ws.Range("B97").Value = 0
$ws.Range("C97").Value = '{''Individualism''}'
$ws.Range("D97").Value = 261
$ws.Range("A98").Value = 'horizon'
$ws.Range("B98").Value = 0
$ws.Range("C98").Value = '{''Tech_optimism''}'
$ws.Range("D98").Value = 99
$ws.Range("A99").Value = 'competition'
$ws.Range("B99").Value = 0
$ws.Range("C99").Value = '{''Whataboutism''}'
$ws.Range("D99").Value = 53
$ws.Range("A100").Value = 'near future'
$ws.Range("B100").Value = 0
$ws.Range("C100").Value = '{''Tech_optimism''}'
$ws.Range("D100").Value = 9
$ws.Range("A101").Value = 'inevitable'
$ws.Range("B101").Value = 0
$ws.Range("C101").Value = '{''Doomism''}'
$ws.Range("D101").Value = 15
$ws.Range("A102").Value = 'invest'
$ws.Range("B102").Value = 0
$ws.Range("C102").Value = '{''Tech_optimism''}'
$ws.Range("D102").Value = 151
$ws.Range("A103").Value = 'rush'
$ws.Range("B103").Value = 0
$ws.Range("C103").Value = '{''Perfect_policy''}'
$ws.Range("D103").Value = 38
$ws.Range("A104").Value = 'adaptation'
$ws.Range("B104").Value = 0
$ws.Range("C104").Value = '{''Change_impossible''}'
$ws.Range("D104").Value = 304
$ws.Range("A105").Value = 'doubt'
$ws.Range("B105").Value = 0
$ws.Range("C105").Value = '{''Change_impossible''}'
$ws.Range("D105").Value = 35
$ws.Range("A106").Value = 'lower carbon'
$ws.Range("B106").Value = 0
$ws.Range("C106").Value = '{''FF_solutionism''}'
$ws.Range("D106").Value = 4
